# Update countries & provincias Spain
# Applies:
#  1. Swap of four pairs of country names (rows that referenced reordered
#     shared-string entries in the original OOXML diff).
#  2. Refreshed case statistics (columns B:H) for a number of country rows.
#  3. Updated "last updated" timestamp string in A1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Timestamp update (A1) -------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 8 de Agosto de 2020 a las 13:05"

# --- 2. Country name swaps (column A) -----------------------------------------
$ws.Range("A54").Value  = "Ghana"
$ws.Range("A55").Value  = "Armenia"

$ws.Range("A151").Value = "Malta"
$ws.Range("A152").Value = "Togo"

$ws.Range("A174").Value = "Islas Feroe"
$ws.Range("A175").Value = "Mongolia"

$ws.Range("A202").Value = "Santa Lucia"
$ws.Range("A203").Value = "Timor Oriental"

# --- 3. Updated statistics (columns B:H) --------------------------------------
# row => @(B, C, D, E, F, G, H)
$updates = @{
    4   = @(5095903, 379,  2618019, 2313762, 0, 28,  164122)
    14  = @(324692,  2125, 282122,  24306,   0, 132, 18264)
    43  = @(62300,   239,  56245,   5699,    0, 0,   356)
    44  = @(60623,   1350, 29872,   28092,   0, 43,  2659)
    54  = @(40533,   436,  37702,   2625,    0, 0,   206)
    55  = @(40185,   200,  32395,   7005,    0, 8,   785)
    57  = @(37054,   39,   25960,   9782,    0, 5,   1312)
    58  = @(36451,   182,  31900,   2565,    0, 0,   1986)
    68  = @(22592,   378,  16313,   6206,    0, 3,   73)
    72  = @(20698,   426,  11320,   9100,    0, 12,  278)
    85  = @(11003,   116,  7329,    3445,    0, 4,   229)
    142 = @(1267,    13,   1115,    146,     0, 0,   6)
    151 = @(1035,    40,   675,     351,     0, 0,   9)
    152 = @(1028,    0,    710,     296,     0, 0,   22)
    174 = @(295,     4,    193,     102,     0, 0,   0)
    175 = @(293,     0,    260,     33,      0, 0,   0)
    181 = @(197,     7,    184,     13,      0, 0,   0)
    183 = @(177,     1,    86,      74,      0, 1,   17)
}

$cols = @("B", "C", "D", "E", "F", "G", "H")
foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range("$($cols[$i])$row").Value = $vals[$i]
    }
}
